$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("db Sizing")

# Update the DB size growth calculation formula in C7
$ws.Range("C7").Formula = "=(IF(ROUNDUP(`$C`$8/25000,0)<=1,1,ROUNDUP(`$C`$8/25000,0))*250)*1024"

# Recalculate the workbook so dependent cells (F7, G7, J8, F9, G9, J9, F10, G10,
# F13, G13, F15, G15, F16, G16, J18, J20, J21, J22, ...) pick up the new values
$excel.CalculateFullRebuild()

# Update the active selection on the sheet to match the saved cursor position
$ws.Range("I28").Select()
